$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("T2").Value = 10
$ws.Range("AA2").Value = 6.5
$ws.Range("AF2").Value = 17
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 4.33
$ws.Range("I6").Value = 7.5
$ws.Range("K6").Value = 15
$ws.Range("R6").Value = 1.73
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 9
$ws.Range("V6").Value = 8.5
$ws.Range("AA6").Value = 9
$ws.Range("G8").Value = 1.5
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 7.5
$ws.Range("K8").Value = 12
$ws.Range("AG8").Value = 19
$ws.Range("I9").Value = 3.75
$ws.Range("K9").Value = 6.1
$ws.Range("L9").Value = 1.36
$ws.Range("M9").Value = 2.9
$ws.Range("Q9").Value = 2.65
$ws.Range("R9").Value = 1.72
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 7.1
$ws.Range("U9").Value = 11
$ws.Range("Y9").Value = 27
$ws.Range("Z9").Value = 6.1
$ws.Range("AA9").Value = 5.3
$ws.Range("AH9").Value = 65
$ws.Range("G10").Value = 1.7
$ws.Range("H10").Value = 3.55
$ws.Range("I10").Value = 4.85
$ws.Range("J10").Value = 1.07
$ws.Range("K10").Value = 6.7
$ws.Range("L10").Value = 1.36
$ws.Range("M10").Value = 2.9
$ws.Range("N10").Value = 2.05
$ws.Range("O10").Value = 1.7
$ws.Range("P10").Value = 1.42
$ws.Range("Q10").Value = 2.67
$ws.Range("R10").Value = 2
$ws.Range("S10").Value = 1.72
$ws.Range("T10").Value = 5.9
$ws.Range("U10").Value = 7.1
$ws.Range("V10").Value = 8.25
$ws.Range("W10").Value = 12.5
$ws.Range("X10").Value = 15
$ws.Range("Z10").Value = 6.7
$ws.Range("AA10").Value = 6.8
$ws.Range("AB10").Value = 18
$ws.Range("AC10").Value = 100
$ws.Range("AE10").Value = 11.5
$ws.Range("AF10").Value = 26
$ws.Range("AG10").Value = 16
$ws.Range("AH10").Value = 90
$ws.Range("AI10").Value = 50
$ws.Range("AJ10").Value = 60
$ws.Range("H11").Value = 3.6
$ws.Range("I11").Value = 4.1
$ws.Range("J11").Value = 1.06
$ws.Range("K11").Value = 10
$ws.Range("L11").Value = 1.36
$ws.Range("M11").Value = 3
$ws.Range("R11").Value = 1.91
$ws.Range("S11").Value = 1.8
$ws.Range("X11").Value = 17
$ws.Range("AH11").Value = 41
$ws.Range("G12").Value = 3.4
$ws.Range("I12").Value = 2
$ws.Range("AB12").Value = 19
$ws.Range("K13").Value = 13
$ws.Range("N13").Value = 1.83
$ws.Range("O13").Value = 1.98
$ws.Range("N14").Value = 1.8
$ws.Range("O14").Value = 2
$ws.Range("L15").Value = 1.29
$ws.Range("M15").Value = 3.5
$ws.Range("N15").Value = 1.98
$ws.Range("O15").Value = 1.83
$ws.Range("G24").Value = 2.63
$ws.Range("H24").Value = 3.25
$ws.Range("I24").Value = 2.7
$ws.Range("J24").Value = 1.07
$ws.Range("K24").Value = 9
$ws.Range("R24").Value = 1.8
$ws.Range("S24").Value = 1.95
$ws.Range("J31").Value = 1.01
$ws.Range("K31").Value = 13
$ws.Range("H32").Value = 3.7
$ws.Range("I32").Value = 1.9
$ws.Range("P32").Value = 1.3
$ws.Range("Q32").Value = 3.4
$ws.Range("R32").Value = 1.62
$ws.Range("S32").Value = 2.2
$ws.Range("Z32").Value = 15
$ws.Range("AA32").Value = 7.5
$ws.Range("AE32").Value = 9.5
$ws.Range("AJ32").Value = 21
$ws.Range("L33").Value = 1.07
$ws.Range("J40").Value = 1.03
$ws.Range("K40").Value = 15
$ws.Range("N40").Value = 1.67
$ws.Range("O40").Value = 2.15
$ws.Range("L42").Value = 1.25
$ws.Range("M42").Value = 3.75
$ws.Range("N42").Value = 1.8
$ws.Range("O42").Value = 2
$ws.Range("L43").Value = 1.25
$ws.Range("M43").Value = 3.75
$ws.Range("N43").Value = 1.88
$ws.Range("O43").Value = 1.93
$ws.Range("G46").Value = 8.5
$ws.Range("H46").Value = 6.5
$ws.Range("I46").Value = 1.22
$ws.Range("N46").Value = 1.2
$ws.Range("O46").Value = 4.33
$ws.Range("P46").Value = 1.13
$ws.Range("Q46").Value = 6
$ws.Range("R46").Value = 1.5
$ws.Range("S46").Value = 2.5
$ws.Range("AA46").Value = 17
$ws.Range("AB46").Value = 19
$ws.Range("AC46").Value = 41
$ws.Range("AD46").Value = 101
$ws.Range("AE46").Value = 19
